# Apply "fix search and toc style" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix shared-string text: search scope description no longer mentions
# excluding only xml/json files ---
$ws.Cells.Item(19, 7).Value = "修改搜索范围,以头信息的“search: exclude”区分是否在搜索范围内"

# --- Add new review rows (22-24) describing search page style issues ---

# Row 22: 王杰 found that the standalone search page style doesn't match
# the overall style; fixed by modifying the search page style.
$ws.Cells.Item(22, 3).Value = "王杰"
$ws.Cells.Item(22, 4).Value = "独立搜索页面样式和整体样式不符"
$ws.Cells.Item(22, 7).Value = "修改搜索页面样式"

# Row 23: header search popup covered by the table of contents.
$ws.Cells.Item(23, 4).Value = "头部搜索弹出框被目录遮挡"

# Row 24: right-hand side table of contents height doesn't line up with text.
$ws.Cells.Item(24, 4).Value = "右侧目录高度未与文本平齐"

# --- Update the view's selection / scroll position to match the new state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 3
$ws.Range("D24").Select()
